# Auto-generated Word COM-interop script to apply regression-table value updates
# per commit: update tables and code for AWM_Hom_grid
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "WARNING: not found -> $old"
    }
}

Replace-Text "-4.88727" "-5.12633"
Replace-Text "0.72093" "0.67657"
Replace-Text "-6.779" "-7.577"
Replace-Text "1.21e-11 ***" "3.54e-14 ***"
Replace-Text "0.29783" "0.30026"
Replace-Text "0.02491" "0.02474"
Replace-Text "11.957" "12.135"
Replace-Text "0.47464" "0.48157"
Replace-Text "0.04108" "0.04013"
Replace-Text "11.555" "12.000"
Replace-Text "0.13744" "0.13135"
Replace-Text "0.04717" "0.04555"
Replace-Text "2.914" "2.884"
Replace-Text "0.003569 **" "0.00393 **"
Replace-Text "-0.38712" "-0.39885"
Replace-Text "0.05536" "0.05468"
Replace-Text "-6.992" "-7.294"
Replace-Text "2.71e-12 ***" "3.02e-13 ***"
Replace-Text "-0.81556" "-0.77893"
Replace-Text "0.08688" "0.08967"
Replace-Text "-9.387" "-8.686"
Replace-Text "0.22987" "0.21442"
Replace-Text "0.06675" "0.06578"
Replace-Text "3.444" "3.260"
Replace-Text "0.000573 ***" "0.00112 **"
Replace-Text "-0.13161" "-0.13513"
Replace-Text "0.05019" "0.05010"
Replace-Text "-2.622" "-2.697"
Replace-Text "0.008739 **" "0.00699 **"
Replace-Text "0.02552" "0.02460"
Replace-Text "0.540" "0.521"
Replace-Text "0.589097" "0.60253"
Replace-Text "AWM_Zip_Hom" "AWM_Hom_grid"
Replace-Text "-0.13329" "-0.07041"
Replace-Text "0.10958" "0.03783"
Replace-Text "-1.216" "-1.861"
Replace-Text "0.223840" "0.06272 ."
Replace-Text "-5.20951" "-2.87028"
Replace-Text "5.50631" "5.22689"
Replace-Text "-0.946" "-0.549"
Replace-Text "0.3441" "0.583"
Replace-Text "0.06590" "0.06616"
Replace-Text "0.19304" "0.19413"
Replace-Text "0.7328" "0.733"
Replace-Text "0.09653" "0.02505"
Replace-Text "0.35806" "0.36169"
Replace-Text "0.270" "0.069"
Replace-Text "0.7875" "0.945"
Replace-Text "0.38935" "0.56083"
Replace-Text "0.44823" "0.43292"
Replace-Text "0.869" "1.295"
Replace-Text "0.3850" "0.195"
Replace-Text "0.68697" "0.67567"
Replace-Text "0.44666" "0.44681"
Replace-Text "1.538" "1.512"
Replace-Text "0.1241" "0.130"
Replace-Text "0.98485" "0.91942"
Replace-Text "0.55802" "0.56176"
Replace-Text "1.765" "1.637"
Replace-Text "0.0776 ." "0.102"
Replace-Text "0.06430" "0.17216"
Replace-Text "0.68618" "0.67971"
Replace-Text "0.094" "0.253"
Replace-Text "0.9253" "0.800"
Replace-Text "-0.14279" "-0.16968"
Replace-Text "0.40057" "0.40025"
Replace-Text "-0.356" "-0.424"
Replace-Text "0.7215" "0.672"
Replace-Text "-0.62736" "-0.61975"
Replace-Text "0.38480" "0.38792"
Replace-Text "-1.630" "-1.598"
Replace-Text "0.1030" "0.110"
Replace-Text "1.46376" "0.15501"
Replace-Text "0.86999" "0.28927"
Replace-Text "1.682" "0.536"
Replace-Text "0.0925 ." "0.592"

Write-Host "Done applying replacements."
